$wb = $excel.ActiveWorkbook

# --- 1. Remove the old "emre" sheet and the old "jan" sheet -------------
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("emre").Delete() | Out-Null
$thomas = $wb.Worksheets.Item("thomas")
$wb.Worksheets.Item("jan").Delete() | Out-Null

# --- 2. Re-create "jan" (same data) right after "thomas" -----------------
$dummy1 = $wb.Worksheets.Add($null, $thomas)
$dummy2 = $wb.Worksheets.Add($null, $dummy1)
$dummy2.Delete() | Out-Null

$newJan = $wb.Worksheets.Add($null, $dummy1)
$newJan.Name = "jan"
$newJan.Range("A1").Value = "groupComponent"
$newJan.Range("B1").Value = "nameComponent"
$newJan.Range("A2").Value = "Motherboard"
$newJan.Range("B2").Value = "GA-Z87-HD3"
$newJan.Range("A3").Value = "RAM"
$newJan.Range("B3").Value = "CML8GX3M2A1600C9"
$newJan.Range("A4").Value = "GPU"
$newJan.Range("B4").Value = "GTX 980"
$newJan.Range("A5").Value = "CPU"
$newJan.Range("B5").Value = "i7 4770K"
$newJan.Range("A6").Value = "PSU"
$newJan.Range("B6").Value = "GS800"
$newJan.Range("A7").Value = "Drive"
$newJan.Range("B7").Value = "950 EVO"

$dummy1.Delete() | Out-Null

# --- 3. Add a brand new "emre" sheet at the end of the tab strip --------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$newEmre = $wb.Worksheets.Add($null, $last)
$newEmre.Name = "emre"

# --- 4. Populate the new "emre" sheet with the (variable) component rows
$newEmre.Range("A1").Value = "groupComponent"
$newEmre.Range("B1").Value = "nameComponent"
$newEmre.Range("A2").Value = "Motherboard"
$newEmre.Range("B2").Value = "X99 Rampage V Extreme"
$newEmre.Range("A3").Value = "RAM"
$newEmre.Range("B3").Value = "HyperX Fury"
$newEmre.Range("A4").Value = "PSU"
$newEmre.Range("B4").Value = "Supernova G2"
$newEmre.Range("A5").Value = "GPU"
$newEmre.Range("B5").Value = "Fury X"

# --- 5. Keep "notworking" as the active tab ------------------------------
$wb.Worksheets.Item("notworking").Activate()
